# Add a new "Play ID" column at the front of the Updated_Play_Database sheet,
# shifting all existing columns (A:L -> B:M) one to the right, and fill the
# new column with sequential Play ID values (P001 .. P013).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A; this shifts every existing
# column (and their widths / contents) one position to the right.
$ws.Range("A1").EntireColumn.Insert()

# Header for the new column.
$ws.Range("A1").Value = "Play ID"

# Sequential Play ID values for the 13 data rows.
$playIds = @("P001","P002","P003","P004","P005","P006","P007","P008","P009","P010","P011","P012","P013")
for ($i = 0; $i -lt $playIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $playIds[$i]
}
